$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlTop = -4160

# Column R (2021) is a new year column added to the right of Q (2020).
# Reuse the neighbouring Q-column cell formatting so the engine resolves to
# the very same (already existing) style indexes instead of minting new ones.

# R3: empty separator cell on the thick-bottom-border row - same look as Q3.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial($xlPasteFormats)

# R4: new "2021" column header - same look as Q4 (2020).
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial($xlPasteFormats)
$ws.Range("R4").Value = 2021

# R5: new sanitary-chemical value - based on Q5's format, but top-aligned
# (this is a genuinely new style in the workbook, vertical="top" instead of
# "center").
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial($xlPasteFormats)
$ws.Range("R5").Value = 0.9
$ws.Range("R5").VerticalAlignment = $xlTop

# R6: new microbiological value - same look as Q6.
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial($xlPasteFormats)
$ws.Range("R6").Value = 6.5

$excel.CutCopyMode = $false

# Update the active selection to reflect where the editor left off.
$null = $ws.Range("T5").Select()
